$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Duplicate the whole "日期" block (rows 222-230, header + table +
#        summary row) down to rows 232-240 to create the entry for the next
#        day (2017.11.15), preserving all formatting/merged cells. Do this
#        before filling in column C below, so the new block keeps blank
#        "完成情况" cells just like the source did at copy time.
$src = $ws.Range("A222:D230")
$dst = $ws.Range("A232:D240")
$src.Copy($dst)

# New block's rows are all short single-line text, so force the normal
# 22.5pt row height (the copied block had some 45pt rows for wrapped text).
$ws.Range("A232:D240").RowHeight = 22.5

# --- 2. Adjust the copied content for the new day's block.
$ws.Range("A232").Value2 = "日期：2017.11.15 第十二周 周三"
$ws.Range("B236").Value2 = "迟到"
$ws.Range("B237").Value2 = "完善前端“我的”模块"
$ws.Range("B239").Value2 = "迟到"
$ws.Range("A240").Value2 = "总结："

# --- 3. Fill in the "完成情况" (completion status) column for the
#        2017.11.13/11.14 block (rows 224-229), which was left blank before.
$ws.Range("C224").Value2 = "未完成"
$ws.Range("C225").Value2 = "未完成"
$ws.Range("C226").Value2 = "完成"
$ws.Range("C227").Value2 = "未完成"
$ws.Range("C228").Value2 = "未完成"
$ws.Range("C229").Value2 = "完成"

# --- 4. Update the summary for the previous (2017.11.13/11.14) block.
$ws.Range("A230").Value2 = "总结：进度缓慢"

# --- 5. Update the selection to reflect where the user last clicked.
$ws.Range("C235").Select()
